$d = $word.ActiveDocument

# 1) "USE A MINER:" -> "VERWENDUNG EINES MINERS:"
$rng = $d.Content
$rng.Find.Execute("USE A MINER:") | Out-Null
$rng.Text = "VERWENDUNG EINES MINERS:"

# 2) "GPU NVIDIA miner" -> "GPU NVIDIA Miner"
$rng = $d.Content
$rng.Find.Execute("GPU NVIDIA miner") | Out-Null
$rng.Text = "GPU NVIDIA Miner"

# 3) "GPU AMD miner" -> "GPU AMD Miner"
$rng = $d.Content
$rng.Find.Execute("GPU AMD miner") | Out-Null
$rng.Text = "GPU AMD Miner"

# 4) "For mining support please join" -> "Für Support hinsichtlich Mining bitte bei " (trailing space)
$rng = $d.Content
$rng.Find.Execute("For mining support please join") | Out-Null
$rng.Text = "Für Support hinsichtlich Mining bitte bei "

# 5) "Discord" -> "Discord melden"
#    This run is preceded by another run in the same paragraph; replacing its
#    whole text via Range.Text re-derives formatting from the insertion point
#    and loses the hyperlink run's rPr. Keep the first character intact so the
#    original run (and its formatting) survives, then adjust the remaining text.
$rng = $d.Content
$rng.Find.Execute("Discord") | Out-Null
$tailStart = $rng.Start + 1
$tailEnd = $rng.End
$tail = $d.Range($tailStart, $tailEnd)
$tail.Text = "iscord melden"

# 6) "<nbsp>EXCHANGES" -> "BÖRSEN" (the original text has a leading non-breaking space)
$rng = $d.Content
$rng.Find.Execute([char]0x00A0 + "EXCHANGES") | Out-Null
$rng.Text = "BÖRSEN"

# 7) "Turn altcoins into Smartcash instantly" -> "Altcoins sofort in Smartcash umtauschen"
$rng = $d.Content
$rng.Find.Execute("Turn altcoins into Smartcash instantly") | Out-Null
$rng.Text = "Altcoins sofort in Smartcash umtauschen"
